$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 90
$ws.Range("C2").Value = 60

$ws.Range("B3").Value = 7

$ws.Range("B4").Value = 18

$ws.Range("B5").Value = 35
$ws.Range("C5").Value = 20

$ws.Range("B6").Value = 50
$ws.Range("C6").Value = 30

$ws.Range("B7").Value = 30
$ws.Range("C7").Value = 20

$ws.Range("B8").Value = 20
$ws.Range("C8").Value = 10

$ws.Range("B9").Value = 12

$ws.Range("B10").Value = 20

$ws.Range("B11").Value = 25
